$d = $word.ActiveDocument

# --- 1) "problem solving behavior" -> "idea generation" -------------------
$rng = $d.Content
$found1 = $rng.Find.Execute(
    " problem solving behavior.", $false, $false, $false, $false, $false,
    $true, 1, $false, " idea generation.", 2)
Write-Output "Replace1: $found1"

# --- 2) "In this project I will be required to ... surveys." --------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "In this project I will be required to work on puzzle and sorting tasks, watch short videos, and complete several surveys. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "In this project I will be required to work on an idea generation task and complete a short survey. I may also work on a categorization task or watch short videos.",
    2)
Write-Output "Replace2: $found2"

# --- 3) Footer version bump: 1.2 -> 1.3 ------------------------------------
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$found3 = $ftr.Range.Find.Execute(
    "Version: 1.2", $false, $false, $false, $false, $false,
    $true, 1, $false, "Version: 1.3", 2)
Write-Output "Replace3: $found3"
